$wb = $excel.ActiveWorkbook

$wsCap = $wb.Worksheets.Item("capabilitiesDetails")
$wsTest = $wb.Worksheets.Item("testData")

# --- capabilitiesDetails sheet ---
# Remove the leading placeholder row (numeric 0..5); the remaining rows
# shift up so the header-label row becomes row 1 and the value row becomes
# row 2.
[void]$wsCap.Rows.Item(1).Delete()

# Update the sheet selection to match the target state (no longer the
# active tab, new selected cell).
[void]$wsCap.Range("A4").Select()

# --- testData sheet ---
# Fix the mis-spelled labels while they are still on row 2 (pre-shift).
$wsTest.Range("B2").Value = "CustomerPassword"
$wsTest.Range("C2").Value = "InvalidPassword"

# Drop the existing hyperlinks; they will be recreated at their shifted
# locations once the placeholder row above them is removed.
[void]$wsTest.Hyperlinks.Delete()

# Remove the leading placeholder row (numeric 0..6); the remaining rows
# shift up.
[void]$wsTest.Rows.Item(1).Delete()

# Recreate the hyperlinks at their new (shifted up) locations.
[void]$wsTest.Hyperlinks.Add($wsTest.Range("A2"), "mailto:user123@gmail.com")
[void]$wsTest.Hyperlinks.Add($wsTest.Range("B2"), "mailto:test@123")
[void]$wsTest.Hyperlinks.Add($wsTest.Range("G2"), "mailto:asd@gk.com")

# Restore the original "Hyperlink" cell style on the linked cells (Add()
# re-applies its own formatting xf).
$wsTest.Range("A2").Style = "Hyperlink"
$wsTest.Range("B2").Style = "Hyperlink"
$wsTest.Range("G2").Style = "Hyperlink"

# testData becomes the active/selected sheet with a new selection.
[void]$wsTest.Select()
[void]$wsTest.Range("C1").Select()
